$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "SACRU - URLs" - 10 URL rows, column A autosized, default (A1) selection
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add($null, $lastSheet)
$ws1.Name = "SACRU - URLs"
$ws1.Range("A1").Value = 'https://www.delta.com/us/en/travel-update-center/extending-skymiles-benefits?src=benefits3'
$ws1.Range("A2").Value = 'http://amex.co/extrasupport?src=amex3'
$ws1.Range("A3").Value = 'https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/onboard-services?src=obsspace1#socialdist'
$ws1.Range("A4").Value = 'https://www.delta.com/us/en/skymiles/program-resources/news-and-updates?src=sm1a#recent'
$ws1.Range("A5").Value = 'https://www.delta.com/us/en/travel-update-center/extending-skymiles-benefits?src=benefits1'
$ws1.Range("A6").Value = 'https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/onboard-services?src=obs1'
$ws1.Range("A7").Value = 'https://www.delta.com/us/en/travel-update-center/flying-what-you-need-to-know/coronavirus-regional-restrictions?src=restrictions1'
$ws1.Range("A8").Value = 'https://www.delta.com/us/en/travel-update-center/from-delta-to-our-customers/medical-volunteers-book-free?src=flyfree1'
$ws1.Range("A9").Value = 'https://www.delta.com/us/en/delta-digital/mobile?src=app1'
$ws1.Range("A10").Value = 'https://www.delta.com/us/en/coronavirus-update-center/ways-we-are-keeping-you-safe/skyclub-update?src=skyclub1'
$ws1.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 2) "OUATT - Header Names" - 3 header rows, portrait page setup, sel C5
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OUATT - Header Names"
# Values entered out of row order so new shared strings land at the indices
# (51, 52, 53) recorded in the workbook: DOWNLOAD... first, then FLYING...,
# then DELTA SKY CLUB UPDATES, finally laid out A1/A2/A3 as below.
$ws2.Range("A2").Value = 'DOWNLOAD THE FLY DELTA APP'
$ws2.Range("A3").Value = 'FLYING WITH A PARTNER AIRLINE?'
$ws2.Range("A1").Value = 'DELTA SKY CLUB UPDATES'
$ws2.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws2.PageSetup.Orientation = 1
$ws2.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) "OUATT - DeltaSkyClub URL" - single URL row, default (A1) selection
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "OUATT - DeltaSkyClub URL"
$ws3.Range("A1").Value = 'https://www.delta.com/us/en/travel-update-center/ways-we-are-keeping-you-safe/skyclub-update'
$ws3.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 4) "OUATT - FlyDeltaApp URL" - single URL row, sel F24
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "OUATT - FlyDeltaApp URL"
$ws4.Range("A1").Value = 'https://www.delta.com/us/en/delta-digital/mobile?src=app2'
$ws4.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws4.Range("F24").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) "OUATT - FlyingPartnerAir URL" - single URL row, no custom width, default selection
#    (sheetId 19 is consumed by a transient sheet removed before the final
#    save, matching the sheetId="20" gap recorded in the workbook)
# ---------------------------------------------------------------------------
$wsTemp = $wb.Worksheets.Add($null, $ws4)
$wsTemp.Name = "TempPlaceholder"

$ws5 = $wb.Worksheets.Add($null, $wsTemp)
$ws5.Name = "OUATT - FlyingPartnerAir URL"
$ws5.Range("A1").Value = 'https://www.delta.com/us/en/travel-update-center/flying-what-you-need-to-know/things-to-know-when-you-travel-with-a-partner-airline'

$excel.DisplayAlerts = $false
$wsTemp.Delete()

# ---------------------------------------------------------------------------
# 6) "FAQ Header Names" - 3 header rows, sel G3 (ends up the active tab)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "FAQ Header Names"
$ws6.Range("A1").Value = 'TRAVEL FLEXIBILITY'
$ws6.Range("A2").Value = 'SAFER TRAVEL'
$ws6.Range("A3").Value = 'DELTA PARTNERS AND AFFILIATES'
$ws6.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws6.Range("G3").Select() | Out-Null
